# Updates the timing-results workbook:
#  - "Brute Force" and "Divide and Conquer" sheets: refresh trial timings for
#    input sizes 1, 50, 100, 200, 400, 800 (rows 2-7) and append new rows for
#    input sizes 1500, 2000, 3000 (rows 8-10), each with AVERAGE/STDEV.S formulas.
#  - "Summary" sheet: mirror the new input sizes and extend the summary
#    formulas down through row 10.
$wb = $excel.ActiveWorkbook

# --- Brute Force sheet ---
$ws = $wb.Worksheets.Item("Brute Force")

$arrA = New-Object "object[,]" 9,1
$arrA[0,0] = 1
$arrA[1,0] = 50
$arrA[2,0] = 100
$arrA[3,0] = 200
$arrA[4,0] = 400
$arrA[5,0] = 800
$arrA[6,0] = 1500
$arrA[7,0] = 2000
$arrA[8,0] = 3000
$ws.Range("A2:A10").Value = $arrA

$arrBK = New-Object "object[,]" 9,10
$arrBK[0,0] = 0.000003576278686523438
$arrBK[0,1] = 0.00001478195190429688
$arrBK[0,2] = 0.000001192092895507812
$arrBK[0,3] = 0.00000095367431640625
$arrBK[0,4] = 0.000001192092895507812
$arrBK[0,5] = 0.00000095367431640625
$arrBK[0,6] = 0.0000007152557373046875
$arrBK[0,7] = 0.0000007152557373046875
$arrBK[0,8] = 0.000002622604370117188
$arrBK[0,9] = 0.000002384185791015625
$arrBK[1,0] = 0.002595901489257812
$arrBK[1,1] = 0.002419233322143555
$arrBK[1,2] = 0.002533435821533203
$arrBK[1,3] = 0.002760887145996094
$arrBK[1,4] = 0.002789497375488281
$arrBK[1,5] = 0.002676010131835938
$arrBK[1,6] = 0.002419471740722656
$arrBK[1,7] = 0.002909660339355469
$arrBK[1,8] = 0.002492189407348633
$arrBK[1,9] = 0.002468585968017578
$arrBK[2,0] = 0.01061725616455078
$arrBK[2,1] = 0.01631331443786621
$arrBK[2,2] = 0.00970458984375
$arrBK[2,3] = 0.009839057922363281
$arrBK[2,4] = 0.009581089019775391
$arrBK[2,5] = 0.01037478446960449
$arrBK[2,6] = 0.01018381118774414
$arrBK[2,7] = 0.009801387786865234
$arrBK[2,8] = 0.01181125640869141
$arrBK[2,9] = 0.009929895401000977
$arrBK[3,0] = 0.04105210304260254
$arrBK[3,1] = 0.03879523277282715
$arrBK[3,2] = 0.03829789161682129
$arrBK[3,3] = 0.04174900054931641
$arrBK[3,4] = 0.03918337821960449
$arrBK[3,5] = 0.03902482986450195
$arrBK[3,6] = 0.04598879814147949
$arrBK[3,7] = 0.0411827564239502
$arrBK[3,8] = 0.04018115997314453
$arrBK[3,9] = 0.03961968421936035
$arrBK[4,0] = 0.1662888526916504
$arrBK[4,1] = 0.1552770137786865
$arrBK[4,2] = 0.1603591442108154
$arrBK[4,3] = 0.1553127765655518
$arrBK[4,4] = 0.1591179370880127
$arrBK[4,5] = 0.1521799564361572
$arrBK[4,6] = 0.1610879898071289
$arrBK[4,7] = 0.1527466773986816
$arrBK[4,8] = 0.1621439456939697
$arrBK[4,9] = 0.1549711227416992
$arrBK[5,0] = 0.6237304210662842
$arrBK[5,1] = 0.6273961067199707
$arrBK[5,2] = 0.6390266418457031
$arrBK[5,3] = 0.6265377998352051
$arrBK[5,4] = 0.6328692436218262
$arrBK[5,5] = 0.6343848705291748
$arrBK[5,6] = 0.6299924850463867
$arrBK[5,7] = 0.6288118362426758
$arrBK[5,8] = 0.6270408630371094
$arrBK[5,9] = 0.6495654582977295
$arrBK[6,0] = 2.225206613540649
$arrBK[6,1] = 2.229511737823486
$arrBK[6,2] = 2.23875880241394
$arrBK[6,3] = 2.249152660369873
$arrBK[6,4] = 2.250405311584473
$arrBK[6,5] = 2.228482961654663
$arrBK[6,6] = 2.313826560974121
$arrBK[6,7] = 2.272916316986084
$arrBK[6,8] = 2.275491237640381
$arrBK[6,9] = 2.26696252822876
$arrBK[7,0] = 3.975231170654297
$arrBK[7,1] = 3.972332715988159
$arrBK[7,2] = 3.961799144744873
$arrBK[7,3] = 3.969098567962646
$arrBK[7,4] = 3.973695039749146
$arrBK[7,5] = 3.962504386901855
$arrBK[7,6] = 3.982734441757202
$arrBK[7,7] = 4.005683660507202
$arrBK[7,8] = 3.913738965988159
$arrBK[7,9] = 3.968324184417725
$arrBK[8,0] = 8.802835464477539
$arrBK[8,1] = 8.894636154174805
$arrBK[8,2] = 8.744202852249146
$arrBK[8,3] = 8.808440923690796
$arrBK[8,4] = 8.718207359313965
$arrBK[8,5] = 8.659956455230713
$arrBK[8,6] = 8.921866178512573
$arrBK[8,7] = 8.795485019683838
$arrBK[8,8] = 9.055953502655029
$arrBK[8,9] = 8.869661808013916
$ws.Range("B2:K10").Value = $arrBK

$ws.Range("L8").Formula = "=AVERAGE(B8:K8)"
$ws.Range("M8").Formula = "=STDEV.S(B8:K8)"
$ws.Range("L9").Formula = "=AVERAGE(B9:K9)"
$ws.Range("M9").Formula = "=STDEV.S(B9:K9)"
$ws.Range("L10").Formula = "=AVERAGE(B10:K10)"
$ws.Range("M10").Formula = "=STDEV.S(B10:K10)"

# --- Divide and Conquer sheet ---
$ws = $wb.Worksheets.Item("Divide and Conquer")

$arrA = New-Object "object[,]" 9,1
$arrA[0,0] = 1
$arrA[1,0] = 50
$arrA[2,0] = 100
$arrA[3,0] = 200
$arrA[4,0] = 400
$arrA[5,0] = 800
$arrA[6,0] = 1500
$arrA[7,0] = 2000
$arrA[8,0] = 3000
$ws.Range("A2:A10").Value = $arrA

$arrBK = New-Object "object[,]" 9,10
$arrBK[0,0] = 0.00003337860107421875
$arrBK[0,1] = 0.0000324249267578125
$arrBK[0,2] = 0.000006437301635742188
$arrBK[0,3] = 0.000005483627319335938
$arrBK[0,4] = 0.000005483627319335938
$arrBK[0,5] = 0.000005245208740234375
$arrBK[0,6] = 0.000004529953002929688
$arrBK[0,7] = 0.000004529953002929688
$arrBK[0,8] = 0.00002455711364746094
$arrBK[0,9] = 0.00001478195190429688
$arrBK[1,0] = 0.0003829002380371094
$arrBK[1,1] = 0.0003759860992431641
$arrBK[1,2] = 0.0003404617309570312
$arrBK[1,3] = 0.0003790855407714844
$arrBK[1,4] = 0.0004532337188720703
$arrBK[1,5] = 0.000457763671875
$arrBK[1,6] = 0.000453948974609375
$arrBK[1,7] = 0.000308990478515625
$arrBK[1,8] = 0.0003244876861572266
$arrBK[1,9] = 0.0003771781921386719
$arrBK[2,0] = 0.0006780624389648438
$arrBK[2,1] = 0.0007014274597167969
$arrBK[2,2] = 0.0006437301635742188
$arrBK[2,3] = 0.0007040500640869141
$arrBK[2,4] = 0.0007100105285644531
$arrBK[2,5] = 0.000789642333984375
$arrBK[2,6] = 0.0006978511810302734
$arrBK[2,7] = 0.0007014274597167969
$arrBK[2,8] = 0.0006926059722900391
$arrBK[2,9] = 0.0006880760192871094
$arrBK[3,0] = 0.001419305801391602
$arrBK[3,1] = 0.001399993896484375
$arrBK[3,2] = 0.001376867294311523
$arrBK[3,3] = 0.001477241516113281
$arrBK[3,4] = 0.00134587287902832
$arrBK[3,5] = 0.001365900039672852
$arrBK[3,6] = 0.001585721969604492
$arrBK[3,7] = 0.001375675201416016
$arrBK[3,8] = 0.001415491104125977
$arrBK[3,9] = 0.001445531845092773
$arrBK[4,0] = 0.002962589263916016
$arrBK[4,1] = 0.003128767013549805
$arrBK[4,2] = 0.003047466278076172
$arrBK[4,3] = 0.002901315689086914
$arrBK[4,4] = 0.002819299697875977
$arrBK[4,5] = 0.002986431121826172
$arrBK[4,6] = 0.002939224243164062
$arrBK[4,7] = 0.002996444702148438
$arrBK[4,8] = 0.003011941909790039
$arrBK[4,9] = 0.002924442291259766
$arrBK[5,0] = 0.006684064865112305
$arrBK[5,1] = 0.007087945938110352
$arrBK[5,2] = 0.006383895874023438
$arrBK[5,3] = 0.006432294845581055
$arrBK[5,4] = 0.01158452033996582
$arrBK[5,5] = 0.00627589225769043
$arrBK[5,6] = 0.009986639022827148
$arrBK[5,7] = 0.006769418716430664
$arrBK[5,8] = 0.006601333618164062
$arrBK[5,9] = 0.00649714469909668
$arrBK[6,0] = 0.01241207122802734
$arrBK[6,1] = 0.01237297058105469
$arrBK[6,2] = 0.01250886917114258
$arrBK[6,3] = 0.01244688034057617
$arrBK[6,4] = 0.0126035213470459
$arrBK[6,5] = 0.01221227645874023
$arrBK[6,6] = 0.01272058486938477
$arrBK[6,7] = 0.01233100891113281
$arrBK[6,8] = 0.01273536682128906
$arrBK[6,9] = 0.01286482810974121
$arrBK[7,0] = 0.01806426048278809
$arrBK[7,1] = 0.01550960540771484
$arrBK[7,2] = 0.01892495155334473
$arrBK[7,3] = 0.01992130279541016
$arrBK[7,4] = 0.01701784133911133
$arrBK[7,5] = 0.0198369026184082
$arrBK[7,6] = 0.01735615730285645
$arrBK[7,7] = 0.01679754257202148
$arrBK[7,8] = 0.01778006553649902
$arrBK[7,9] = 0.01632046699523926
$arrBK[8,0] = 0.02908515930175781
$arrBK[8,1] = 0.02730035781860352
$arrBK[8,2] = 0.02657079696655273
$arrBK[8,3] = 0.02612066268920898
$arrBK[8,4] = 0.02603816986083984
$arrBK[8,5] = 0.02605795860290527
$arrBK[8,6] = 0.02683424949645996
$arrBK[8,7] = 0.02833366394042969
$arrBK[8,8] = 0.02717375755310059
$arrBK[8,9] = 0.02778482437133789
$ws.Range("B2:K10").Value = $arrBK

$ws.Range("L8").Formula = "=AVERAGE(B8:K8)"
$ws.Range("M8").Formula = "=STDEV.S(B8:K8)"
$ws.Range("L9").Formula = "=AVERAGE(B9:K9)"
$ws.Range("M9").Formula = "=STDEV.S(B9:K9)"
$ws.Range("L10").Formula = "=AVERAGE(B10:K10)"
$ws.Range("M10").Formula = "=STDEV.S(B10:K10)"

# --- Summary sheet ---
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("A3").Value = 50
$ws.Range("A4").Value = 100
$ws.Range("A5").Value = 200
$ws.Range("A6").Value = 400
$ws.Range("A7").Value = 800
$ws.Range("A8").Value = 1500
$ws.Range("A9").Value = 2000
$ws.Range("A10").Value = 3000

$ws.Range("B8").Formula = "='Brute Force'!L8"
$ws.Range("C8").Formula = "='Brute Force'!M8"
$ws.Range("D8").Formula = "='Divide and Conquer'!L8"
$ws.Range("E8").Formula = "='Divide and Conquer'!M8"
$ws.Range("F8").Formula = "=B8/D8"

$ws.Range("B9").Formula = "='Brute Force'!L9"
$ws.Range("C9").Formula = "='Brute Force'!M9"
$ws.Range("D9").Formula = "='Divide and Conquer'!L9"
$ws.Range("E9").Formula = "='Divide and Conquer'!M9"
$ws.Range("F9").Formula = "=B9/D9"

$ws.Range("B10").Formula = "='Brute Force'!L10"
$ws.Range("C10").Formula = "='Brute Force'!M10"
$ws.Range("D10").Formula = "='Divide and Conquer'!L10"
$ws.Range("E10").Formula = "='Divide and Conquer'!M10"
$ws.Range("F10").Formula = "=B10/D10"

